# Jogos_da_Semana_FlashScore_2024-11-04.xlsx update
# - updates a handful of odds on existing rows 2 and 3
# - inserts a brand-new match row at row 8 (COLOMBIA - PRIMERA A,
#   Atl. Nacional vs Santa Fe), pushing the old rows 8-10 down to 9-11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Atl. Tucuman vs Sarmiento Junin) odds tweaks ---
$ws.Range("G2").Value = 1.9
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.63
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("AI2").Value = 15
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 11
$ws.Range("AV2").Value = 67
$ws.Range("AX2").Value = 26
$ws.Range("BD2").Value = 151

# --- Row 3 (Ind. Rivadavia vs Rosario Central) odds tweak ---
$ws.Range("N3").Value = 5

# --- Insert a new row at position 8; existing rows 8-10 shift to 9-11 ---
$ws.Rows(8).Insert()

# New row 8: COLOMBIA - PRIMERA A, Atl. Nacional vs Santa Fe
$ws.Range("A8").Value = "txqKnEdc"
# Force the Date column to stay plain text (it would otherwise be parsed
# as a date serial number, unlike the rest of the sheet's inline strings).
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "04/11/2024"
$ws.Range("C8").Value = "22:15"
$ws.Range("D8").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E8").Value = "Atl. Nacional"
$ws.Range("F8").Value = "Santa Fe"
$ws.Range("G8").Value = 1.8
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 2.6
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 5.5
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.44
$ws.Range("P8").Value = 2.63
$ws.Range("Q8").Value = 2.5
$ws.Range("R8").Value = 1.5
$ws.Range("S8").Value = 1.53
$ws.Range("T8").Value = 2.38
$ws.Range("U8").Value = 2.25
$ws.Range("V8").Value = 1.57
$ws.Range("W8").Value = 5.5
$ws.Range("X8").Value = 7.5
$ws.Range("Y8").Value = 9.5
$ws.Range("Z8").Value = 15
$ws.Range("AA8").Value = 19
$ws.Range("AB8").Value = 41
$ws.Range("AC8").Value = 6.5
$ws.Range("AD8").Value = 6.5
$ws.Range("AE8").Value = 21
$ws.Range("AF8").Value = 81
$ws.Range("AG8").Value = 10
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 17
$ws.Range("AJ8").Value = 51
$ws.Range("AK8").Value = 41
$ws.Range("AL8").Value = 51
$ws.Range("AM8").Value = 900
$ws.Range("AN8").Value = 3.6
$ws.Range("AO8").Value = 10
$ws.Range("AP8").Value = 26
$ws.Range("AQ8").Value = 41
$ws.Range("AR8").Value = 67
$ws.Range("AS8").Value = 251
$ws.Range("AT8").Value = 2.38
$ws.Range("AU8").Value = 9.5
$ws.Range("AV8").Value = 81
$ws.Range("AW8").Value = 6.5
$ws.Range("AX8").Value = 29
$ws.Range("AY8").Value = 41
$ws.Range("AZ8").Value = 101
$ws.Range("BA8").Value = 151
$ws.Range("BB8").Value = 351
$ws.Range("BC8").Value = 126
$ws.Range("BD8").Value = 126
